$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.068707790979985361
$ws.Range("A2").Value = -0.046061934238363023
$ws.Range("A3").Value = -0.0089999995131453403
$ws.Range("A4").Value = 0.28399186991552128
$ws.Range("A5").Value = -0.0059999995277983942
$ws.Range("A6").Value = -0.0059999995128947603
$ws.Range("A7").Value = -0.019999999432442905
$ws.Range("A8").Value = -0.019999999429767712
$ws.Range("A9").Value = -0.005999999507188214
$ws.Range("A10").Value = -0.0059999995059172306
$ws.Range("A11").Value = -0.0044999995143584215
$ws.Range("A12").Value = -0.0059999995057307132
$ws.Range("A13").Value = -0.0059999995062609557
$ws.Range("A14").Value = -0.011999999472727119
$ws.Range("A15").Value = 0.042853778580768775
$ws.Range("A16").Value = -0.0059999995065060929
$ws.Range("A17").Value = -0.0059999995046924326
$ws.Range("A18").Value = -0.0089999994873064537
$ws.Range("A19").Value = -0.010864574018345952
$ws.Range("A20").Value = -0.0089999995083935858
$ws.Range("A21").Value = -0.0089999995078153816
$ws.Range("A22").Value = -0.0089999995074254713
$ws.Range("A23").Value = -0.0089999995091059048
$ws.Range("A24").Value = -0.04199999931408005
$ws.Range("A25").Value = -0.041999999310557534
$ws.Range("A26").Value = -0.005999999511359988
$ws.Range("A27").Value = -0.0059999995090720404
$ws.Range("A28").Value = -0.0059999994994361927
$ws.Range("A29").Value = -0.011999999459016308
$ws.Range("A30").Value = -0.01999999941077979
$ws.Range("A31").Value = -0.013664947504221203
$ws.Range("A32").Value = -0.020999999399953673
$ws.Range("A33").Value = -0.0059999994848576321
